$d = $word.ActiveDocument

# --------------------------------------------------------------------
# 1) Merge the "is the degrees [gramStart]latitude.[gramEnd] " runs
#    (with the stray grammar-check proofErr bookmarks) into a single
#    plain run reading " is the degrees latitude. ". Find/Replace
#    across the run boundary causes Word to collapse the runs and
#    drop the now-irrelevant proofErr markers.
# --------------------------------------------------------------------
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute(
    " is the degrees latitude. ",  # FindText
    $false,                        # MatchCase
    $false,                        # MatchWholeWord
    $false,                        # MatchWildcards
    $false,                        # MatchSoundsLike
    $false,                        # MatchAllWordForms
    $true,                         # Forward
    1,                             # Wrap (wdFindContinue)
    $false,                        # Format
    " is the degrees latitude. ",  # ReplaceWith
    2                              # Replace (wdReplaceAll)
) | Out-Null

# --------------------------------------------------------------------
# 2) Append two new paragraphs at the end of the document: a blank
#    paragraph, and one holding the new sensitivity-tests sentence
#    (Times font, matching the rest of the body text).
# --------------------------------------------------------------------
$lastPara = $d.Paragraphs.Last
$lastPara.Range.InsertParagraphAfter() | Out-Null

$blankPara = $d.Paragraphs.Last
$blankPara.Range.InsertParagraphAfter() | Out-Null

$newPara = $d.Paragraphs.Last
$newPara.Range.Text = "We also conduct a suite of sensitivity tests to provide additional constraints on the error of the optimized emissions, which are summarized in section 2.6."

# Re-fetch the paragraph (Text= may reseat ranges) and apply the Times
# font to the run's text only, leaving the trailing paragraph mark
# untouched so no stray pPr/rPr is introduced on the paragraph mark.
$newPara = $d.Paragraphs.Last
$textRange = $newPara.Range
$body = $d.Range($textRange.Start, $textRange.End - 1)
$body.Font.NameAscii = "Times"
$body.Font.Name = "Times"
